# Apply data corrections for 4-25-2014-15 team stats sheet.
# The underlying NBA stats were captured one calendar day off, so several
# ranking/rate columns need to shift and the embedded date label needs to be
# normalized from "4-25-2014-15" to the ISO-style "2015-04-25".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric corrections ---
$cellUpdates = @{
    "AS2" = 22
    "AH3" = 8
    "AQ3" = 14
    "AY3" = 21
    "AZ3" = 20
    "AH4" = 3
    "AT4" = 22
    "AU4" = 20
    "AR5" = 26
    "AY5" = 23
    "AF6" = 9
    "AG6" = 9
    "AH6" = 8
    "AV6" = 12
    "AH7" = 23
    "AP7" = 12
    "AQ7" = 18
    "AS7" = 21
    "AA8" = 22.1
    "AB8" = 105.2
    "AC8" = 2.9
    "AD8" = 1
    "AF8" = 9
    "AJ8" = 8
    "AQ8" = 16
    "AS8" = 23
    "AT8" = 23
    "AX8" = 20
    "D8" = 82
    "F8" = 32
    "G8" = 0.61
    "J8" = 85.8
    "K8" = 0.463
    "N8" = 0.352
    "Q8" = 0.752
    "S8" = 31.8
    "T8" = 42.3
    "U8" = 22.5
    "AH9" = 8
    "AS9" = 12
    "AX9" = 21
    "AJ10" = 6
    "BA10" = 24
    "BB10" = 18
    "AA12" = 21.1
    "AC12" = 3.4
    "AD12" = 1
    "AE12" = 3
    "AG12" = 3
    "AO12" = 5
    "AP12" = 2
    "AT12" = 14
    "AY12" = 22
    "D12" = 82
    "E12" = 56
    "G12" = 0.6830000000000001
    "J12" = 83.3
    "N12" = 0.348
    "O12" = 18.6
    "P12" = 26
    "Q12" = 0.715
    "T12" = 43.7
    "Y12" = 5.3
    "Z12" = 22
    "AQ13" = 13
    "AA15" = 19.4
    "AB15" = 98.5
    "AC15" = -6.8
    "AD15" = 1
    "AP15" = 11
    "AS15" = 13
    "AT15" = 12
    "AU15" = 21
    "AX15" = 22
    "AZ15" = 21
    "BA15" = 23
    "BB15" = 19
    "D15" = 82
    "F15" = 61
    "G15" = 0.256
    "I15" = 37.2
    "J15" = 85.59999999999999
    "K15" = 0.435
    "M15" = 18.9
    "O15" = 17.5
    "P15" = 23.6
    "Q15" = 0.741
    "S15" = 32.3
    "T15" = 43.9
    "U15" = 20.9
    "Y15" = 4.8
    "Z15" = 21.2
    "AE16" = 5
    "AK16" = 9
    "AL17" = 21
    "AX17" = 18
    "AH18" = 3
    "AH19" = 20
    "AP19" = 3
    "AA20" = 18.7
    "AB20" = 99.40000000000001
    "AC20" = 0.8
    "AD20" = 1
    "AE20" = 13
    "AG20" = 13
    "AI20" = 11
    "AJ20" = 19
    "AK20" = 10
    "AQ20" = 17
    "AS20" = 19
    "AW20" = 25
    "BC20" = 13
    "D20" = 82
    "E20" = 45
    "G20" = 0.549
    "J20" = 82.90000000000001
    "K20" = 0.457
    "L20" = 7.1
    "N20" = 0.37
    "P20" = 21.8
    "Q20" = 0.751
    "R20" = 11.5
    "S20" = 32
    "T20" = 43.5
    "X20" = 6.2
    "AO22" = 6
    "AQ22" = 15
    "AB23" = 95.7
    "AC23" = -5.7
    "AD23" = 1
    "AJ23" = 20
    "AL23" = 22
    "AR23" = 25
    "AY23" = 24
    "AZ23" = 19
    "D23" = 82
    "F23" = 57
    "G23" = 0.305
    "I23" = 37.5
    "K23" = 0.453
    "M23" = 19.5
    "S23" = 31.8
    "T23" = 41.8
    "U23" = 20.6
    "Y23" = 5.4
    "AH24" = 20
    "AJ24" = 22
    "AH25" = 8
    "AJ25" = 6
    "AX26" = 18
    "AW27" = 26
    "AC28" = 6.2
    "AD28" = 1
    "AE28" = 5
    "AG28" = 5
    "AP28" = 23
    "AT28" = 15
    "AV28" = 13
    "D28" = 82
    "E28" = 55
    "G28" = 0.671
    "J28" = 83.59999999999999
    "K28" = 0.468
    "L28" = 8.300000000000001
    "N28" = 0.367
    "P28" = 21.4
    "S28" = 33.8
    "T28" = 43.6
    "V28" = 14
    "AZ29" = 18
    "AP30" = 13
    "AS30" = 20
    "AT30" = 11
    "AH31" = 1
    "AJ31" = 21
    "AP31" = 22
    "BB31" = 17
    "BC31" = 14
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}

# --- Date label correction (keep stored as text, not a date serial) ---
$dateRows = 2..31
$ws.Range("BF2:BF31").NumberFormat = "@"
foreach ($r in $dateRows) {
    $ws.Cells.Item($r, 58).Value = "2015-04-25"
}
